# Updates mechanism labels to be more descriptive, and moves the active
# selection to G23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the abbreviated "Mechanism" values (column B) into fuller names.
$mechanismMap = @{
    "HMGCoA"            = "HMGCoA inhibitor"
    "RAAS"              = "RAAS inhibitor"
    "ECM"               = "ECM modulator"
    "IC Enz Inhibitor"  = "Intracellular Enzyme Inhibitor"
    "GF Modulator"      = "Growth Factor Modulator"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    if ($current -ne $null -and $mechanismMap.ContainsKey($current)) {
        $cell.Value2 = $mechanismMap[$current]
    }
}

# Move the selection as recorded in the saved view state.
$ws.Range("G23").Select()
